$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing contents, since the new data occupies a smaller range
# than the original sheet (A1:I26 -> A1:I12).
$ws.Cells.Clear()

# --- Section 1 header (Level) ---
$ws.Range("A1").Value = "Level"

# Column headers for section 1
$ws.Range("A2").Value = "Rank"
$ws.Range("B2").Value = "Title"
$ws.Range("C2").Value = "Creator"
$ws.Range("D2").Value = "Date Posted"
$ws.Range("E2").Value = "Country"
$ws.Range("F2").Value = "Language"
$ws.Range("G2").Value = "TUS (2024-04-22)"
$ws.Range("H2").Value = "Rating"
$ws.Range("I2").Value = "Comment Count"

# Section 1 data rows (rows 3-6)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "诀别书"
$ws.Range("C3").Value = "可乐没有气了"
$ws.Range("D3").Value = "19/03/2024"
$ws.Range("E3").Value = "N/A"
$ws.Range("F3").Value = "zh-cn"
$ws.Range("G3").Value = 184
$ws.Range("H3").Value = "N/A"
$ws.Range("I3").Value = "0"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "BOB BALL3"
$ws.Range("C4").Value = "脸红"
$ws.Range("D4").Value = "21/03/2024"
$ws.Range("E4").Value = "CN"
$ws.Range("F4").Value = "zh-cn"
$ws.Range("G4").Value = 1518
$ws.Range("H4").Value = "N/A"
$ws.Range("I4").Value = "0"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "禁止摆烂萌新图"
$ws.Range("C5").Value = "肉女士"
$ws.Range("D5").Value = "23/03/2024"
$ws.Range("E5").Value = "N/A"
$ws.Range("F5").Value = "zh-cn"
$ws.Range("G5").Value = 1195
$ws.Range("H5").Value = "N/A"
$ws.Range("I5").Value = "0"

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "PartyTest_Fyang"
$ws.Range("C6").Value = "肥羊"
$ws.Range("D6").Value = "27/03/2024"
$ws.Range("E6").Value = "CN"
$ws.Range("F6").Value = "zh-cn"
$ws.Range("G6").Value = 500
$ws.Range("H6").Value = "N/A"
$ws.Range("I6").Value = "0"

# --- Section 2 header (Model) ---
$ws.Range("A7").Value = "Model"

# Column headers for section 2
$ws.Range("A8").Value = "Rank"
$ws.Range("B8").Value = "Title"
$ws.Range("C8").Value = "Creator"
$ws.Range("D8").Value = "Date Posted"
$ws.Range("E8").Value = "Country"
$ws.Range("F8").Value = "Language"
$ws.Range("G8").Value = "TUS (2024-04-22)"
$ws.Range("H8").Value = "Rating"
$ws.Range("I8").Value = "Comment Count"

# Section 2 data rows (rows 9-12)
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Dogday"
$ws.Range("C9").Value = "Kimmel"
$ws.Range("D9").Value = "23/03/2024"
$ws.Range("E9").Value = "N/A"
$ws.Range("F9").Value = "es"
$ws.Range("G9").Value = 35
$ws.Range("H9").Value = "N/A"
$ws.Range("I9").Value = "0"

$ws.Range("A10").Value = 2
$ws.Range("B10").Value = "Adolf Hitler"
$ws.Range("C10").Value = "ebrunedre"
$ws.Range("D10").Value = "23/03/2024"
$ws.Range("E10").Value = "TR"
$ws.Range("F10").Value = "tr"
$ws.Range("G10").Value = 456
$ws.Range("H10").Value = "N/A"
$ws.Range("I10").Value = "2"

$ws.Range("A11").Value = 3
$ws.Range("B11").Value = "RED ROBIN!"
$ws.Range("C11").Value = "ordinalst"
$ws.Range("D11").Value = "07/03/2024"
$ws.Range("E11").Value = "N/A"
$ws.Range("F11").Value = "en"
$ws.Range("G11").Value = 117
$ws.Range("H11").Value = "N/A"
$ws.Range("I11").Value = "0"

$ws.Range("A12").Value = 4
$ws.Range("B12").Value = "КЛОУН"
$ws.Range("C12").Value = "7700n"
$ws.Range("D12").Value = "02/03/2024"
$ws.Range("E12").Value = "FR"
$ws.Range("F12").Value = "ru"
$ws.Range("G12").Value = 1053
$ws.Range("H12").Value = "N/A"
$ws.Range("I12").Value = "0"

# Some values above look like numbers/dates to Excel's auto-detection
# (the Comment Count "0"/"2" values, and the Date Posted "07/03/2024" /
# "02/03/2024" strings, which parse as valid M/D/Y dates since day<=12).
# Re-enter those specific cells while temporarily flagged as Text so the
# literal string is preserved, then clear the format override again so
# the cell keeps its default (General) style, matching a plain text
# shared-string cell.
$textFixups = @{
  "I3"  = "0"
  "I4"  = "0"
  "I5"  = "0"
  "I6"  = "0"
  "I9"  = "0"
  "I10" = "2"
  "I11" = "0"
  "I12" = "0"
  "D11" = "07/03/2024"
  "D12" = "02/03/2024"
}
foreach ($addr in $textFixups.Keys) {
  $cell = $ws.Range($addr)
  $cell.NumberFormat = "@"
  $cell.Value = $textFixups[$addr]
  $cell.ClearFormats()
}
